$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update pin numbers that are now documented with specific AVR port names
$ws.Range("B21").Value = "27-PC4"
$ws.Range("B22").Value = "28-PC5"

# Add newly documented connections for pins 4 (INT0) and 5 (INT1)
$ws.Range("C10").Value = "Solonoid "
$ws.Range("C9").Value = "LCD-!RST"

# Update the active selection to reflect where the edit was made
$ws.Range("C10").Select()
